$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)
$ws4 = $wb.Worksheets.Item(4)
$ws5 = $wb.Worksheets.Item(5)

# Rename sheets
$ws1.Name = "GNG_TO-16509961231865957"
$ws2.Name = "NB_TO-16509961253865628"
$ws3.Name = "RS_TO-16509961253865628"
$ws4.Name = "TOL_TO-16509961254585967"
$ws5.Name = "vSAT_TO-16509961255385635"

# Sheet1 (GNG_TO) updates
$ws1.Range("B2").Value = "go_stims-16509961231545618.csv"
$ws1.Range("B3").Value = "GNG_stims-1650996123170602.csv"
$ws1.Range("B4").Value = "go_stims-1650996123170602.csv"
$ws1.Range("B5").Value = "GNG_stims-16509961231865957.csv"

# Sheet2 (NB_TO) updates
$ws2.Range("B2").Value = "OB-16509961241385598.csv"
$ws2.Range("B3").Value = "ZB-match_6-165099612336257.csv"
$ws2.Range("B4").Value = "OB-16509961249946468.csv"
$ws2.Range("B5").Value = "ZB-match_5-16509961238825977.csv"
$ws2.Range("B6").Value = "TB-1650996125218605.csv"
$ws2.Range("B7").Value = "OB-16509961243866117.csv"
$ws2.Range("B8").Value = "TB-16509961253706124.csv"
$ws2.Range("B9").Value = "ZB-match_1-16509961235065846.csv"
$ws2.Range("B10").Value = "TB-16509961250346012.csv"

# Sheet4 (TOL_TO) updates
$ws4.Range("B2").Value = "MM_stims-16509961254185662.csv"
$ws4.Range("B3").Value = "ZM_stims-16509961253945613.csv"
$ws4.Range("B4").Value = "MM_stims-16509961254425964.csv"
$ws4.Range("B5").Value = "ZM_stims-16509961254185662.csv"
$ws4.Range("B6").Value = "MM_stims-16509961254585967.csv"
$ws4.Range("B7").Value = "ZM_stims-16509961254425964.csv"

# Sheet5 (vSAT_TO) updates
$ws5.Range("B2").Value = "SAT_stims-1650996125490617.csv"
$ws5.Range("B3").Value = "vSAT_stims-16509961255065856.csv"
$ws5.Range("B4").Value = "vSAT_stims-16509961255225677.csv"
$ws5.Range("B5").Value = "SAT_stims-16509961254665756.csv"
